# Add a new "Course" column in front of the existing exam-data table.
#
# The author didn't use Excel's "Insert Column" command (that would also
# have shifted every <col> width definition one slot to the right); instead
# the existing columns A-H were retyped/shifted one cell to the right by
# hand and a new value was entered in column A. So here we shift each row's
# values rightward (starting from the rightmost column so nothing is
# clobbered before it's read) and then fill in the new first column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 3
$lastCol = 8   # original table used columns A:H

for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = ($lastCol + 1); $c -ge 2; $c--) {
        $src = $ws.Cells.Item($r, $c - 1)
        $dst = $ws.Cells.Item($r, $c)
        $dst.Value = $src.Value2
    }
}

# New leading "Course" column with its per-row codes.
$ws.Range("A1").Value = "Course"
$ws.Range("A2").Value = "A1"
$ws.Range("A3").Value = "B1"

# The column itself got narrower, and the two rightmost columns (now H & I)
# were widened/added by hand after the retype.
$ws.Columns("A:A").ColumnWidth = 18.83
$ws.Columns("H:H").ColumnWidth = 33
$ws.Columns("I:I").ColumnWidth = 27.5

# Leave the selection where the author finished editing.
$ws.Range("F11").Select() | Out-Null
